$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = [double]"9.01614080345281e-05"
$ws.Range("H2").Value = 0.003556098334722702
$ws.Range("K2").Value = 5.347171369788641
$ws.Range("L2").Value = "[2.2221450152771, 8.472197724300182]"
$ws.Range("M2").Value = 0.0008451594919505201
$ws.Range("N2").Value = 0.00169031898390104
$ws.Range("O2").Value = -1.056631763369539
$ws.Range("P2").Value = "[-1.6730002920017704, -0.44026323473730855]"
$ws.Range("Q2").Value = 0.0008266124890827786
$ws.Range("R2").Value = 0.0008266124890827786
$ws.Range("S2").Value = 14.27974321893345
$ws.Range("T2").Value = "[12.5687136688529, 15.990772769014]"
$ws.Range("W2").Value = 4.370690690690797
$ws.Range("X2").Value = 1.821121121121169
$ws.Range("Y2").Value = 6.920260260260424

# Row 3 updates
$ws.Range("E3").Value = 24.07000000000032
$ws.Range("G3").Value = 0.0004790421237154119
$ws.Range("H3").Value = 0.003556098334722702
$ws.Range("K3").Value = 5.289868836112849
$ws.Range("L3").Value = "[1.733016479841993, 8.846721192383704]"
$ws.Range("M3").Value = 0.003694183817008589
$ws.Range("N3").Value = 0.003694183817008589
$ws.Range("O3").Value = 2.182447749340657
$ws.Range("P3").Value = "[1.540921321580579, 2.8239741771007347]"
$ws.Range("Q3").Value = [double]"1.137028249331706e-10"
$ws.Range("R3").Value = [double]"2.274056498663413e-10"
$ws.Range("S3").Value = 12.71008832269407
$ws.Range("T3").Value = "[10.853099495016584, 14.567077150371547]"
$ws.Range("W3").Value = 15.70934934934956
$ws.Range("X3").Value = 13.25175175175193
$ws.Range("Y3").Value = 18.16694694694719
